$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 533.3333
$ws.Range("I2").Value = 533.3333
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 533.3333
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -420.3333
$ws.Range("N2").ClearContents()

$ws.Range("H15").Value = 975.717
$ws.Range("I15").Value = 975.717
$ws.Range("K15").Value = 2927.151
$ws.Range("M15").Value = -2758.151

$ws.Range("H40").Value = 1739.3
$ws.Range("I40").Value = 2078.6
$ws.Range("J40").Value = 1400
$ws.Range("K40").Value = 2078.6
$ws.Range("L40").Value = 1400
$ws.Range("M40").Value = -1903.6
$ws.Range("N40").Value = -1750

$ws.Range("H43").Value = 1000
$ws.Range("J43").Value = 1000
$ws.Range("L43").Value = 1000
$ws.Range("N43").Value = -1138

$ws.Range("H101").Value = 4041.3333
$ws.Range("I101").Value = 5249.5
$ws.Range("J101").Value = 1625
$ws.Range("K101").Value = 15748.5
$ws.Range("L101").Value = 4875
$ws.Range("M101").Value = -14126.5
$ws.Range("N101").Value = -8119

$ws.Range("H116").Value = 15566.667
$ws.Range("J116").Value = 5157.143
$ws.Range("L116").Value = 5157.143
$ws.Range("N116").Value = -12041.143

$ws.Range("H132").Value = 1504
$ws.Range("I132").Value = 1373.9131
$ws.Range("K132").Value = 4121.7393
$ws.Range("M132").Value = -1591.7393

$ws.Range("H137").Value = 2701.4285
$ws.Range("J137").Value = 2871.9
$ws.Range("L137").Value = 8615.700000000001
$ws.Range("N137").Value = -13715.7

$ws.Range("H138").Value = 3688.25
$ws.Range("I138").Value = 7929
$ws.Range("J138").Value = 2531.682
$ws.Range("K138").Value = 23787
$ws.Range("L138").Value = 7595.045999999999
$ws.Range("M138").Value = -18647
$ws.Range("N138").Value = -17875.046

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4166.635
$ws.Range("I32").Value = 3218.309
$ws.Range("K32").Value = 3218.309
$ws.Range("M32").Value = -2931.309

$ws.Range("H74").Value = 1131.9714
$ws.Range("I74").Value = 565.3214
$ws.Range("K74").Value = 565.3214
$ws.Range("M74").Value = 308.6786

$ws.Range("H77").Value = 1131.9714
$ws.Range("I77").Value = 565.3214
$ws.Range("K77").Value = 2826.607
$ws.Range("M77").Value = 1541.393

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 764.2
$ws.Range("I94").Value = 764.2
$ws.Range("K94").Value = 764.2
$ws.Range("M94").Value = -313.2

$ws.Range("H99").Value = 1298.5
$ws.Range("I99").Value = 949.5
$ws.Range("J99").Value = 1996.5
$ws.Range("K99").Value = 949.5
$ws.Range("L99").Value = 1996.5
$ws.Range("M99").Value = 548.5
$ws.Range("N99").Value = -4992.5

$ws.Range("H134").Value = 5672.1377
$ws.Range("I134").Value = 6230.52
$ws.Range("K134").Value = 18691.56
$ws.Range("M134").Value = -16156.56

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2803.889
$ws.Range("I31").Value = 2873.125
$ws.Range("K31").Value = 2873.125
$ws.Range("M31").Value = -2578.125

$ws.Range("H34").Value = 2803.889
$ws.Range("I34").Value = 2873.125
$ws.Range("K34").Value = 2873.125
$ws.Range("M34").Value = -2671.125

$ws.Range("H99").Value = 2622.9333
$ws.Range("I99").Value = 2032.2727
$ws.Range("J99").Value = 4247.25
$ws.Range("K99").Value = 2032.2727
$ws.Range("L99").Value = 4247.25
$ws.Range("M99").Value = -534.2727
$ws.Range("N99").Value = -7243.25

$ws.Range("H107").Value = 815.6
$ws.Range("I107").Value = 584.0833
$ws.Range("J107").Value = 1741.6666
$ws.Range("K107").Value = 584.0833
$ws.Range("L107").Value = 1741.6666
$ws.Range("M107").Value = 1335.9167
$ws.Range("N107").Value = -5581.6666

$ws.Range("H122").Value = 5346.778
$ws.Range("I122").Value = 4158.4287
$ws.Range("K122").Value = 12475.2861
$ws.Range("M122").Value = -10025.2861

$ws.Range("H126").Value = 2622.9333
$ws.Range("I126").Value = 2032.2727
$ws.Range("J126").Value = 4247.25
$ws.Range("K126").Value = 6096.8181
$ws.Range("L126").Value = 12741.75
$ws.Range("M126").Value = -3626.8181
$ws.Range("N126").Value = -17681.75

$ws.Range("H132").Value = 5379.5
$ws.Range("I132").Value = 4670.3335
$ws.Range("J132").Value = 5805
$ws.Range("K132").Value = 14011.0005
$ws.Range("L132").Value = 17415
$ws.Range("M132").Value = -11481.0005
$ws.Range("N132").Value = -22475

$ws.Range("H134").Value = 3635.1428
$ws.Range("I134").Value = 3210.6667
$ws.Range("K134").Value = 9632.000100000001
$ws.Range("M134").Value = -7097.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 70000
$ws.Range("J37").Value = 70000
$ws.Range("L37").Value = 210000
$ws.Range("N37").Value = -210224

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H131").Value = 16588.205
$ws.Range("J131").Value = 17746.365
$ws.Range("L131").Value = 53239.095
$ws.Range("N131").Value = -63319.095

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4051.2856
$ws.Range("I70").Value = 3000
$ws.Range("J70").Value = 4226.5
$ws.Range("K70").Value = 3000
$ws.Range("L70").Value = 4226.5
$ws.Range("M70").Value = -2730
$ws.Range("N70").Value = -4766.5

$ws.Range("H73").Value = 4051.2856
$ws.Range("I73").Value = 3000
$ws.Range("J73").Value = 4226.5
$ws.Range("K73").Value = 3000
$ws.Range("L73").Value = 4226.5
$ws.Range("M73").Value = -2064
$ws.Range("N73").Value = -6098.5

$ws.Range("H102").Value = 4625.5
$ws.Range("I102").Value = 5560.077
$ws.Range("J102").Value = 3275.5557
$ws.Range("K102").Value = 5560.077
$ws.Range("L102").Value = 3275.5557
$ws.Range("M102").Value = -3938.077
$ws.Range("N102").Value = -6519.5557

$ws.Range("H122").Value = 1525.48
$ws.Range("I122").Value = 1392.4445
$ws.Range("J122").Value = 1867.5714
$ws.Range("K122").Value = 4177.333500000001
$ws.Range("L122").Value = 5602.7142
$ws.Range("M122").Value = -1727.333500000001
$ws.Range("N122").Value = -10502.7142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 15134.818
$ws.Range("I40").Value = 23400.8
$ws.Range("J40").Value = 8246.5
$ws.Range("K40").Value = 23400.8
$ws.Range("L40").Value = 8246.5
$ws.Range("M40").Value = -23264.8
$ws.Range("N40").Value = -8518.5

$ws.Range("H46").Value = 1523.0588
$ws.Range("I46").Value = 1096.125
$ws.Range("J46").Value = 1902.5555
$ws.Range("K46").Value = 1096.125
$ws.Range("L46").Value = 1902.5555
$ws.Range("M46").Value = -908.125
$ws.Range("N46").Value = -2278.5555

$ws.Range("H122").Value = 5724.4165
$ws.Range("I122").Value = 5149.125
$ws.Range("J122").Value = 6875
$ws.Range("K122").Value = 15447.375
$ws.Range("L122").Value = 20625
$ws.Range("M122").Value = -12997.375
$ws.Range("N122").Value = -25525

$ws.Range("H135").Value = 36150
$ws.Range("J135").Value = 36150
$ws.Range("L135").Value = 36150
$ws.Range("N135").Value = -46290

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 12840.1
$ws.Range("I126").Value = 14650.25
$ws.Range("K126").Value = 43950.75
$ws.Range("M126").Value = -41480.75
Write-Host "Edits applied successfully"
